$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (column labels shift: new Jan_2026 column inserted, Oct_2025 dropped)
$ws.Range("A1").Value = "ISIN"
$ws.Range("B1").Value = "Stock Name"
$ws.Range("C1").Value = "Mutual Fund"
$ws.Range("D1").Value = "Jan_2026"
$ws.Range("E1").Value = "Dec_2025"
$ws.Range("F1").Value = "Nov_2025"
$ws.Range("G1").Value = "MoM"
$ws.Range("H1").Value = "QoQ"

# Row 2: Sequent Scientific Limited
$ws.Range("A2").Value = "INE807F01027"
$ws.Range("B2").Value = "Sequent Scientific Limited"
$ws.Range("C2").Value = "quant Healthcare Fund"
$ws.Range("D2").Value = 8.831270999999999
$ws.Range("E2").Value = 8.50142
$ws.Range("F2").Value = 8.719996999999999
$ws.Range("G2").Value = 0.3298509999999997
$ws.Range("H2").Value = 0.1112739999999999

# Row 3: Zydus Wellness Ltd
$ws.Range("A3").Value = "INE768C01028"
$ws.Range("B3").Value = "Zydus Wellness Ltd"
$ws.Range("C3").Value = "quant Healthcare Fund"
$ws.Range("D3").Value = 8.447488999999999
$ws.Range("E3").Value = 8.051306
$ws.Range("F3").Value = 7.297499
$ws.Range("G3").Value = 0.3961829999999988
$ws.Range("H3").Value = 1.149989999999999

# Row 4: SMS Pharmaceuticals Limited
$ws.Range("A4").Value = "INE812G01025"
$ws.Range("B4").Value = "SMS Pharmaceuticals Limited"
$ws.Range("C4").Value = "quant Healthcare Fund"
$ws.Range("D4").Value = 8.441843
$ws.Range("E4").Value = 9.582952000000001
$ws.Range("F4").Value = 9.657165000000001
$ws.Range("G4").Value = -1.141109
$ws.Range("H4").Value = -1.215322

# Row 5: Aurobindo Pharma Limited
$ws.Range("A5").Value = "INE406A01037"
$ws.Range("B5").Value = "Aurobindo Pharma Limited"
$ws.Range("C5").Value = "quant Healthcare Fund"
$ws.Range("D5").Value = 7.705935
$ws.Range("E5").Value = 7.094618
$ws.Range("F5").Value = 7.042017
$ws.Range("G5").Value = 0.6113170000000006
$ws.Range("H5").Value = 0.6639179999999998

# Row 6: Aster DM Healthcare Limited
$ws.Range("A6").Value = "INE914M01019"
$ws.Range("B6").Value = "Aster DM Healthcare Limited"
$ws.Range("C6").Value = "quant Healthcare Fund"
$ws.Range("D6").Value = 7.241881
$ws.Range("E6").Value = 7.604171
$ws.Range("F6").Value = 7.852396
$ws.Range("G6").Value = -0.3622899999999998
$ws.Range("H6").Value = -0.6105149999999995

# Row 7: Alivus Life Sciences
$ws.Range("A7").Value = "INE03Q201024"
$ws.Range("B7").Value = "Alivus Life Sciences"
$ws.Range("C7").Value = "quant Healthcare Fund"
$ws.Range("D7").Value = 6.950458
$ws.Range("E7").Value = 6.230534
$ws.Range("F7").Value = 5.846258
$ws.Range("G7").Value = 0.7199240000000007
$ws.Range("H7").Value = 1.104200000000001

# Row 8: Pfizer Ltd
$ws.Range("A8").Value = "INE182A01018"
$ws.Range("B8").Value = "Pfizer Ltd"
$ws.Range("C8").Value = "quant Healthcare Fund"
$ws.Range("D8").Value = 6.425974
$ws.Range("E8").Value = 6.601
$ws.Range("F8").Value = 6.302559
$ws.Range("G8").Value = -0.1750259999999999
$ws.Range("H8").Value = 0.1234150000000005

# Row 9: ANTHEM BIOSCIENCES LIMITED
$ws.Range("A9").Value = "INE0CZ201020"
$ws.Range("B9").Value = "ANTHEM BIOSCIENCES LIMITED"
$ws.Range("C9").Value = "quant Healthcare Fund"
$ws.Range("D9").Value = 5.869715
$ws.Range("E9").Value = 5.997429
$ws.Range("F9").Value = 5.592856
$ws.Range("G9").Value = -0.1277140000000001
$ws.Range("H9").Value = 0.276859

# Row 10: Adani Enterprises Limited
$ws.Range("A10").Value = "INE423A01024"
$ws.Range("B10").Value = "Adani Enterprises Limited"
$ws.Range("C10").Value = "quant Healthcare Fund"
$ws.Range("D10").Value = 5.635718
$ws.Range("E10").Value = 5.871914
$ws.Range("F10").Value = 5.722387
$ws.Range("G10").Value = -0.2361960000000005
$ws.Range("H10").Value = -0.08666900000000055

# Row 11: Divi's Laboratories Limited
$ws.Range("A11").Value = "INE361B01024"
$ws.Range("B11").Value = "Divi's Laboratories Limited"
$ws.Range("C11").Value = "quant Healthcare Fund"
$ws.Range("D11").Value = 3.628896
$ws.Range("E11").Value = 3.603264
$ws.Range("F11").Value = 3.494731
$ws.Range("G11").Value = 0.02563200000000032
$ws.Range("H11").Value = 0.1341650000000003

# Row 12: Medplus Health Services Limited
$ws.Range("A12").Value = "INE804L01022"
$ws.Range("B12").Value = "Medplus Health Services Limited"
$ws.Range("C12").Value = "quant Healthcare Fund"
$ws.Range("D12").Value = 3.520651
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 3.520651
$ws.Range("H12").Value = 3.520651

# Row 13: Concord Biotech Limited
$ws.Range("A13").Value = "INE338H01029"
$ws.Range("B13").Value = "Concord Biotech Limited"
$ws.Range("C13").Value = "quant Healthcare Fund"
$ws.Range("D13").Value = 2.880442
$ws.Range("E13").Value = 3.099855
$ws.Range("F13").Value = 3.126033
$ws.Range("G13").Value = -0.2194129999999999
$ws.Range("H13").Value = -0.2455910000000001

# Row 14: Gland Pharma Limited
$ws.Range("A14").Value = "INE068V01023"
$ws.Range("B14").Value = "Gland Pharma Limited"
$ws.Range("C14").Value = "quant Healthcare Fund"
$ws.Range("D14").Value = 1.52235
$ws.Range("E14").Value = 1.334298
$ws.Range("F14").Value = 1.305113
$ws.Range("G14").Value = 0.1880520000000001
$ws.Range("H14").Value = 0.2172370000000001

# Row 15: Sun Pharmaceutical Industries Limited
$ws.Range("A15").Value = "INE044A01036"
$ws.Range("B15").Value = "Sun Pharmaceutical Industries Limited"
$ws.Range("C15").Value = "quant Healthcare Fund"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 4.800678
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = -4.800678

# Row 16: Adani Enterprises Limited Rights
$ws.Range("A16").Value = "INE423A20016"
$ws.Range("B16").Value = "Adani Enterprises Limited Rights"
$ws.Range("C16").Value = "quant Healthcare Fund"
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0.133224
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = -0.133224

# Row 17: Godrej Properties Limited
$ws.Range("A17").Value = "INE484J01027"
$ws.Range("B17").Value = "Godrej Properties Limited"
$ws.Range("C17").Value = "quant Healthcare Fund"
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 2.626417
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = -2.626417

# Row 18: Samvardhana Motherson International Ltd
$ws.Range("A18").Value = "INE775A01035"
$ws.Range("B18").Value = "Samvardhana Motherson International Ltd"
$ws.Range("C18").Value = "quant Healthcare Fund"
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 4.091851
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = -4.091851
$ws.Range("H18").Value = 0

# Row 19: Cipla Limited
$ws.Range("A19").Value = "INE059A01026"
$ws.Range("B19").Value = "Cipla Limited"
$ws.Range("C19").Value = "quant Healthcare Fund"
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 10.010286
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = -10.010286
$ws.Range("H19").Value = 0
